# Commit: "Add simple Insert data private all other add simple data"
#
# The diff adds six repeated 3x3 blocks of plain (text) data to Sheet1:
#   - four copies in columns F:H at rows 16-18, 24-26, 32-34, 40-42
#   - two copies in columns A:C at rows 43-45, 46-48
#
# Each 3x3 block looks like:
#   34    2r3   34
#   1     3we   34
#   wer1  3wer  34wr
#
# All of these values must land as TEXT cells (shared strings), even the
# ones that look numeric ("34", "1"). Plain `Range.Value = "34"` would be
# auto-coerced to a number by the engine (like typing into the Excel UI),
# so for values that parse as a number we stage the text in a scratch cell
# (forced to text with a leading apostrophe) and use Copy/PasteSpecial
# (values only) into the destination - the paste carries the text type
# without carrying the scratch cell's quote-prefix style along with it.
# The scratch work happens on a throwaway worksheet that is removed again
# at the end, so it leaves no residue in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$scratchSheet = $wb.Worksheets.Add()
$scratch = $scratchSheet.Cells.Item(1, 1)

function Set-CellValue($ws, $addr, $val) {
    if ($val -match '^[0-9]+(\.[0-9]+)?$') {
        # Numeric-looking text: force text via apostrophe prefix in a
        # scratch cell, then copy/paste-values so the destination keeps
        # the default (unstyled) cell format.
        $scratch.Value = "'" + $val
        $scratch.Copy() | Out-Null
        $ws.Range($addr).PasteSpecial(-4163) | Out-Null
    } else {
        # Not numeric-looking - Excel stores it as text on its own.
        $ws.Range($addr).Value = $val
    }
}

Set-CellValue $ws1 'F16' '34'
Set-CellValue $ws1 'G16' '2r3'
Set-CellValue $ws1 'H16' '34'
Set-CellValue $ws1 'F17' '1'
Set-CellValue $ws1 'G17' '3we'
Set-CellValue $ws1 'H17' '34'
Set-CellValue $ws1 'F18' 'wer1'
Set-CellValue $ws1 'G18' '3wer'
Set-CellValue $ws1 'H18' '34wr'
Set-CellValue $ws1 'F24' '34'
Set-CellValue $ws1 'G24' '2r3'
Set-CellValue $ws1 'H24' '34'
Set-CellValue $ws1 'F25' '1'
Set-CellValue $ws1 'G25' '3we'
Set-CellValue $ws1 'H25' '34'
Set-CellValue $ws1 'F26' 'wer1'
Set-CellValue $ws1 'G26' '3wer'
Set-CellValue $ws1 'H26' '34wr'
Set-CellValue $ws1 'F32' '34'
Set-CellValue $ws1 'G32' '2r3'
Set-CellValue $ws1 'H32' '34'
Set-CellValue $ws1 'F33' '1'
Set-CellValue $ws1 'G33' '3we'
Set-CellValue $ws1 'H33' '34'
Set-CellValue $ws1 'F34' 'wer1'
Set-CellValue $ws1 'G34' '3wer'
Set-CellValue $ws1 'H34' '34wr'
Set-CellValue $ws1 'F40' '34'
Set-CellValue $ws1 'G40' '2r3'
Set-CellValue $ws1 'H40' '34'
Set-CellValue $ws1 'F41' '1'
Set-CellValue $ws1 'G41' '3we'
Set-CellValue $ws1 'H41' '34'
Set-CellValue $ws1 'F42' 'wer1'
Set-CellValue $ws1 'G42' '3wer'
Set-CellValue $ws1 'H42' '34wr'
Set-CellValue $ws1 'A43' '34'
Set-CellValue $ws1 'B43' '2r3'
Set-CellValue $ws1 'C43' '34'
Set-CellValue $ws1 'A44' '1'
Set-CellValue $ws1 'B44' '3we'
Set-CellValue $ws1 'C44' '34'
Set-CellValue $ws1 'A45' 'wer1'
Set-CellValue $ws1 'B45' '3wer'
Set-CellValue $ws1 'C45' '34wr'
Set-CellValue $ws1 'A46' '34'
Set-CellValue $ws1 'B46' '2r3'
Set-CellValue $ws1 'C46' '34'
Set-CellValue $ws1 'A47' '1'
Set-CellValue $ws1 'B47' '3we'
Set-CellValue $ws1 'C47' '34'
Set-CellValue $ws1 'A48' 'wer1'
Set-CellValue $ws1 'B48' '3wer'
Set-CellValue $ws1 'C48' '34wr'

$scratchSheet.Delete() | Out-Null
